$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: ABC / S5 / 10.1 / 1.1
$ws.Range("A5").Value = "ABC"
# Row 6: BBC / S6 / 10.1 / 1.1
$ws.Range("A6").Value = "BBC"
# Row 5/6 status column (creates shared strings S5, S6 in order)
$ws.Range("B5").Value = "S5"
$ws.Range("B6").Value = "S6"
# Row 7: CDC / S3 / (no C) / 3
$ws.Range("A7").Value = "CDC"
$ws.Range("B7").Value = "S3"

$ws.Range("C5").Value = 10.1
$ws.Range("D5").Value = 1.1

$ws.Range("C6").Value = 10.1
$ws.Range("D6").Value = 1.1

$ws.Range("D7").Value = 3

[void]$ws.Range("C7").Select()
